$wb = $excel.ActiveWorkbook

# ============================================================
# 1) Rename sheet "Powerplants2022" -> "Powerplants2020"
# ============================================================
$ws3 = $wb.Worksheets.Item("Powerplants2022")
$ws3.Name = "Powerplants2020"

# ============================================================
# 2) Sheet "Overview": append row 11 (copy of row 9, A11 = 9)
# ============================================================
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("A11").Value = 9
$ws1.Range("B11").Value = 2020
$ws1.Range("C11").Value = 544999999.1
$ws1.Range("D11").Value = 20811313152.32112
$ws1.Range("E11").Value = 38.1858957553916
$ws1.Range("F11").Value = 11
$ws1.Range("G11").Value = 222307.2994627971
$ws1.Range("H11").Value = 0
$ws1.Range("I11").Value = 2.560297453063297
$ws1.Range("J11").Value = 97504.61057238668
$ws1.Range("K11").Value = 17769.20341851865
$ws1.Range("L11").Value = 0.1822396224568984
$ws1.Range("M11").Value = 0
$ws1.Range("N11").Value = 0
$ws1.Range("O11").Value = 0
$ws1.Range("P11").Value = 0
$ws1.Range("Q11").Value = 37.85801221665666
$ws1.Range("A9").Copy()
$ws1.Range("A11").PasteSpecial(-4122)

# ============================================================
# 3) Sheet "Capacity": append row 11 (copy of row 9, A11 = 9)
# ============================================================
$ws2 = $wb.Worksheets.Item("Capacity")
$ws2.Range("A11").Value = 9
$ws2.Range("B11").Value = 2020
$ws2.Range("C11").Value = 4644.4034
$ws2.Range("D11").Value = 25208582.8382924
$ws2.Range("E11").Value = 954932128.661841
$ws2.Range("F11").Value = 37.88123016623044
$ws2.Range("G11").Value = 24845.77
$ws2.Range("H11").Value = 59003.61621933627
$ws2.Range("I11").Value = 3332129.912351787
$ws2.Range("J11").Value = 56.47331682121212
$ws2.Range("K11").Value = 31358.329
$ws2.Range("L11").Value = 217108263.1152519
$ws2.Range("M11").Value = 8460882562.993628
$ws2.Range("N11").Value = 38.97079936797324
$ws2.Range("O11").Value = 8194.3025
$ws2.Range("P11").Value = 1007081.087896536
$ws2.Range("Q11").Value = 53413355.2968122
$ws2.Range("R11").Value = 53.03779004367492
$ws2.Range("S11").Value = 8858.749999999998
$ws2.Range("T11").Value = 18624635.99999999
$ws2.Range("U11").Value = 705091697.2187846
$ws2.Range("V11").Value = 37.85801221665675
$ws2.Range("W11").Value = 8599
$ws2.Range("X11").Value = 74259244.44886312
$ws2.Range("Y11").Value = 2846665877.166388
$ws2.Range("Z11").Value = 38.33416160228613
$ws2.Range("AA11").Value = 47547.50848700004
$ws2.Range("AB11").Value = 80823362.9723005
$ws2.Range("AC11").Value = 2945418492.530766
$ws2.Range("AD11").Value = 36.44266192610929
$ws2.Range("AE11").Value = 10271.8
$ws2.Range("AF11").Value = 40873004.28171189
$ws2.Range("AG11").Value = 1517259559.244995
$ws2.Range("AH11").Value = 37.12131236518557
$ws2.Range("AI11").Value = 53555.51607579708
$ws2.Range("AJ11").Value = 50877740.2720072
$ws2.Range("AK11").Value = 1831210598.874818
$ws2.Range("AL11").Value = 35.9923728743579
$ws2.Range("AM11").Value = 20779.02
$ws2.Range("AN11").Value = 36159080.46745713
$ws2.Range("AO11").Value = 1493106750.420732
$ws2.Range("AP11").Value = 41.29271903815462
$ws2.Range("AQ11").Value = 3652.9
$ws2.Range("AR11").Value = 0
$ws2.Range("AS11").Value = 0
$ws2.Range("AT11").Value = 0
$ws2.Range("AU11").Value = 0
$ws2.Range("AV11").Value = 0
$ws2.Range("AW11").Value = 0
$ws2.Range("AX11").Value = 0
$ws2.Range("AY11").Value = 0
$ws2.Range("AZ11").Value = 0
$ws2.Range("BA11").Value = 0
$ws2.Range("BB11").Value = 0
$ws2.Range("A9").Copy()
$ws2.Range("A11").PasteSpecial(-4122)

# ============================================================
# 4) Sheet "Powerplants2020" (was "Powerplants2022"):
#    delete row 6 (decommissioned in-pipeline CCGT that never
#    got built) and refresh the recalculated figures for the
#    remaining plants.
# ============================================================
$ws3.Rows.Item(6).Delete()

$ws3.Range("G2").Value = 21
$ws3.Range("K2").Value = 195047557.2328703
$ws3.Range("M2").Value = 47896307.39275556
$ws3.Range("N2").Value = 25208582.8382924
$ws3.Range("O2").Value = 954932128.661841
$ws3.Range("P2").Value = 711988264.0362153
$ws3.Range("Q2").Value = 37.88123016623044
$ws3.Range("G3").Value = 29
$ws3.Range("K3").Value = 211791260.5158438
$ws3.Range("M3").Value = 7488169619.082788
$ws3.Range("N3").Value = 217108263.1152519
$ws3.Range("O3").Value = 8460882562.993628
$ws3.Range("P3").Value = 760921683.3949952
$ws3.Range("Q3").Value = 38.97079936797324
$ws3.Range("G4").Value = 15
$ws3.Range("K4").Value = 288176673.3787938
$ws3.Range("M4").Value = 109111540.0126057
$ws3.Range("N4").Value = 80823362.9723005
$ws3.Range("O4").Value = 2945418492.530766
$ws3.Range("P4").Value = 2548130279.139367
$ws3.Range("Q4").Value = 36.44266192610929
$ws3.Range("G5").Value = 32
$ws3.Range("K5").Value = 320825454.4060951
$ws3.Range("M5").Value = 3184620.871803447
$ws3.Range("N5").Value = 59003.61621933627
$ws3.Range("O5").Value = 3332129.912351787
$ws3.Range("P5").Value = -320677945.3655468
$ws3.Range("Q5").Value = 56.47331682121212
$ws3.Range("A6").Value = 4
$ws3.Range("G6").Value = 37
$ws3.Range("K6").Value = 5150718.870612221
$ws3.Range("P6").Value = -5150718.870612221
$ws3.Range("A7").Value = 5
$ws3.Range("G7").Value = 57
$ws3.Range("K7").Value = 7384220.949805131
$ws3.Range("O7").Value = 705091697.2187846
$ws3.Range("P7").Value = 697707476.2689794
$ws3.Range("Q7").Value = 37.85801221665675
$ws3.Range("A8").Value = 6
$ws3.Range("G8").Value = 39
$ws3.Range("K8").Value = 190684909.5157455
$ws3.Range("M8").Value = 1345556085.273859
$ws3.Range("N8").Value = 36159080.46745713
$ws3.Range("O8").Value = 1493106750.420732
$ws3.Range("P8").Value = -43134244.36887294
$ws3.Range("Q8").Value = 41.29271903815462
$ws3.Range("A9").Value = 7
$ws3.Range("G9").Value = 36
$ws3.Range("K9").Value = 165046505.4929781
$ws3.Range("M9").Value = 450056026.9628067
$ws3.Range("N9").Value = 74259244.44886312
$ws3.Range("O9").Value = 2846665877.166388
$ws3.Range("P9").Value = 2231563344.710603
$ws3.Range("Q9").Value = 38.33416160228613
$ws3.Range("A10").Value = 8
$ws3.Range("G10").Value = 29
$ws3.Range("K10").Value = 15418557.38030332
$ws3.Range("M10").Value = 49274830.62694587
$ws3.Range("N10").Value = 1007081.087896536
$ws3.Range("O10").Value = 53413355.2968122
$ws3.Range("P10").Value = -11280032.71043699
$ws3.Range("Q10").Value = 53.03779004367492
$ws3.Range("A11").Value = 9
$ws3.Range("G11").Value = 11
$ws3.Range("K11").Value = 227017807.7274436
$ws3.Range("O11").Value = 1831210598.874818
$ws3.Range("P11").Value = 1604192791.147375
$ws3.Range("Q11").Value = 35.9923728743579
$ws3.Range("A12").Value = 10
$ws3.Range("G12").Value = 7
$ws3.Range("K12").Value = 263186052.2797029
$ws3.Range("M12").Value = 110357111.5606221
$ws3.Range("N12").Value = 40873004.28171189
$ws3.Range("O12").Value = 1517259559.244995
$ws3.Range("P12").Value = 1143716395.40467
$ws3.Range("Q12").Value = 37.12131236518557
